$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.136152267456055
$ws.Range("B1").Value = 2.282155990600586
$ws.Range("C1").Value = 10.34047031402588
$ws.Range("D1").Value = 2.169620990753174
$ws.Range("E1").Value = 1.275270104408264
